# Workbook was edited per commit "modified tables except the ones from
# opinion and buy of....": correct the Rolling Stones row label (it
# referenced the wrong/old Band ID 431, row actually holds Band ID 781),
# underline the ACDC row's id cell (A5), and leave the active selection on
# the merged label cell that was being edited (D4:H4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled band id text in the merged D4:H4 cell (431 -> 781,
# matching A4's actual value of 781).
$ws.Range("D4").Value = "#The Rolling Stones´s Band ID: 781"

# Leave the selection where the edit happened.
$ws.Range("D4:H4").Select()

# Underline the id cell on the ACDC row (A5).
$ws.Range("A5").Font.Underline = $true

# Re-affirm J4's formula (CONCATENATE building the INSERT statement).
$ws.Range("J4").Formula = '=CONCATENATE($N$1,A4,$N$3,B4,$N$2)'
